# Practica 4: Version Final del codigo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duty (E5): 1 -> 0.5
$ws.Range("E5").Value = 0.5

# valPR2/PWM period formula (E7): simplify (Duty*Tpwm)/((1/Fosc)*Pre) -> (Duty*Tpwm*Fosc)/(Pre)
$ws.Range("E7").Formula = "=(Duty*Tpwm*Fosc)/(Pre)"

# Update the saved view state: active cell/selection moves to E7, zoom changes to 205%
$ws.Range("E7").Select() | Out-Null
$excel.ActiveWindow.Zoom = 205
